{"js": "// Update the date title and the 25 three-digit-by-one-digit multiplication\n// answers in the worksheet table, per the target diff.\n\n// 1) Update the title paragraph (date line).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.trim() === \"2025-07-28 Monday\") {\n  titlePara.insertText(\"2025-07-29 Tuesday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the table of answers. The table holds the computed answers in\n// row-major order; empty spacer rows are left untouched. We replace the\n// whole `values` grid in one shot (this keeps each cell's existing run\n// formatting intact, only the text content changes), mapping old -> new\n// answer by its position in the table (NOT by text-matching), since one of\n// the new values (\"960\u00d77=6720\") collides with an old value that sits\n// earlier in the table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst oldToNew = {\n  \"601\u00d75=3005\": \"301\u00d73=903\",\n  \"610\u00d73=1830\": \"505\u00d74=2020\",\n  \"568\u00d77=3976\": \"459\u00d75=2295\",\n  \"945\u00d78=7560\": \"502\u00d77=3514\",\n  \"770\u00d72=1540\": \"695\u00d73=2085\",\n  \"740\u00d78=5920\": \"732\u00d72=1464\",\n  \"762\u00d78=6096\": \"410\u00d78=3280\",\n  \"224\u00d78=1792\": \"591\u00d79=5319\",\n  \"906\u00d73=2718\": \"696\u00d72=1392\",\n  \"138\u00d76=828\": \"206\u00d78=1648\",\n  \"762\u00d76=4572\": \"771\u00d77=5397\",\n  \"960\u00d77=6720\": \"949\u00d73=2847\",\n  \"434\u00d78=3472\": \"213\u00d75=1065\",\n  \"105\u00d74=420\": \"534\u00d74=2136\",\n  \"513\u00d78=4104\": \"506\u00d76=3036\",\n  \"737\u00d78=5896\": \"458\u00d74=1832\",\n  \"346\u00d76=2076\": \"436\u00d75=2180\",\n  \"737\u00d73=2211\": \"702\u00d78=5616\",\n  \"742\u00d78=5936\": \"602\u00d77=4214\",\n  \"725\u00d77=5075\": \"960\u00d77=6720\",\n  \"816\u00d78=6528\": \"701\u00d73=2103\",\n  \"279\u00d78=2232\": \"764\u00d72=1528\",\n  \"701\u00d72=1402\": \"604\u00d73=1812\",\n  \"934\u00d76=5604\": \"292\u00d72=584\",\n  \"652\u00d79=5868\": \"195\u00d74=780\",\n};\n\nconst newValues = table.values.map((row) =>\n  row.map((cell) => (Object.prototype.hasOwnProperty.call(oldToNew, cell) ? oldToNew[cell] : cell))\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the date title and the 25 three-digit-by-one-digit multiplication\n# answers in the worksheet table, per the target diff.\n\n$d = $word.ActiveDocument\n\n# 1) Update the title paragraph (date line).\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text.TrimEnd(\"`r`a`n\") -eq \"2025-07-28 Monday\") {\n    $titlePara.Range.Text = \"2025-07-29 Tuesday\"\n}\n\n# 2) Update the table of answers. The table holds the computed answers in\n# row-major order; empty spacer rows are left untouched. We walk every\n# cell and replace its text by exact old -> new lookup (NOT a global\n# text-search/replace), since one of the new values (\"960\u00d77=6720\")\n# collides with an old value that sits earlier in the table, and a\n# naive find/replace could touch it twice.\n$map = @{\n    \"601\u00d75=3005\" = \"301\u00d73=903\"\n    \"610\u00d73=1830\" = \"505\u00d74=2020\"\n    \"568\u00d77=3976\" = \"459\u00d75=2295\"\n    \"945\u00d78=7560\" = \"502\u00d77=3514\"\n    \"770\u00d72=1540\" = \"695\u00d73=2085\"\n    \"740\u00d78=5920\" = \"732\u00d72=1464\"\n    \"762\u00d78=6096\" = \"410\u00d78=3280\"\n    \"224\u00d78=1792\" = \"591\u00d79=5319\"\n    \"906\u00d73=2718\" = \"696\u00d72=1392\"\n    \"138\u00d76=828\"  = \"206\u00d78=1648\"\n    \"762\u00d76=4572\" = \"771\u00d77=5397\"\n    \"960\u00d77=6720\" = \"949\u00d73=2847\"\n    \"434\u00d78=3472\" = \"213\u00d75=1065\"\n    \"105\u00d74=420\"  = \"534\u00d74=2136\"\n    \"513\u00d78=4104\" = \"506\u00d76=3036\"\n    \"737\u00d78=5896\" = \"458\u00d74=1832\"\n    \"346\u00d76=2076\" = \"436\u00d75=2180\"\n    \"737\u00d73=2211\" = \"702\u00d78=5616\"\n    \"742\u00d78=5936\" = \"602\u00d77=4214\"\n    \"725\u00d77=5075\" = \"960\u00d77=6720\"\n    \"816\u00d78=6528\" = \"701\u00d73=2103\"\n    \"279\u00d78=2232\" = \"764\u00d72=1528\"\n    \"701\u00d72=1402\" = \"604\u00d73=1812\"\n    \"934\u00d76=5604\" = \"292\u00d72=584\"\n    \"652\u00d79=5868\" = \"195\u00d74=780\"\n}\n\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cellText = $cell.Range.Text.TrimEnd(\"`r`a`n\")\n        if ($map.ContainsKey($cellText)) {\n            $cell.Range.Text = $map[$cellText]\n        }\n    }\n}\n"}
